$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-07 (row 20)
$ws.Range("B20").Value = 6153
$ws.Range("D20").Value = 5563286
$ws.Range("E20").Value = 904.1582967658053
$ws.Range("F20").Value = 6.287787182587667
$ws.Range("H20").Value = 25.8470373176211
